$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven cell updates: "A1Ref|NewValue" pairs taken from the refreshed
# cryptos (coinranking.com) snapshot - price (D) + 1h change (E) columns,
# plus three ranking-order swaps that moved rows: Mantle<->Kaspa (38/39) and
# Monero/FLOKI/VeChain (48/49/50), each carrying the coin name + link too.
$updates = @'
D2|68.773.53
E2|  -0.51%  
D3|3.868.19
E3|  -1.42%  
E4|  -0.06%  
D5|602.11
E5|  -0.29%  
D6|168.68
E6|  +2.42%  
D7|3.868.27
E7|  -1.44%  
E8|  -0.07%  
D9|0.531
E9|  -0.33%  
D10|0.166
E10|  -0.91%  
D11|6.36
E11|  -0.41%  
D12|0.465
E12|  +0.34%  
D13|37.75
E13|  +1.65%  
D14|0.0000251
E14|  +2.00%  
D15|4.514.72
E15|  -1.41%  
D16|3.855.57
E16|  -3.06%  
D17|68.884.33
E17|  -0.47%  
D18|18.51
E18|  +7.72%  
D19|7.58
E19|  +0.92%  
E20|  -1.12%  
D21|10.85
E21|  -3.71%  
D22|479.47
E22|  -2.05%  
D23|0.740
E23|  +2.12%  
D24|0.0000162
E24|  -3.04%  
D25|84.78
E25|  +0.32%  
D26|2.26
E26|  -0.40%  
D27|12.37
E27|  +1.79%  
D28|10.11
E28|  -0.09%  
E29|  -0.03%  
D30|2.97
E30|  +0.75%  
D31|4.016.82
E31|  -1.34%  
D32|7.79
E32|  -1.37%  
D33|2.32
E33|  -2.96%  
D34|31.30
E34|  -3.92%  
D35|3.834.60
E35|  -0.80%  
D36|0.106
E36|  -1.34%  
D37|6.01
E37|  +1.49%  
B38|Mantle
C38|https://coinranking.com/coin/BoI4ux0nd+mantle-mnt
D38|1.02
E38|  -2.03%  
B39|Kaspa
C39|https://coinranking.com/coin/V8GxkwWow+kaspa-kas
D39|0.140
E39|  -0.30%  
D40|3.30
E40|  +7.58%  
D41|0.999
E41|  -0.09%  
D42|0.318
E42|  -0.72%  
D43|2.04
E43|  +2.02%  
D44|430.19
E44|  -2.45%  
D45|47.86
E45|  -1.21%  
D47|8.64
E47|  +1.96%  
B48|Monero
C48|https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr
D48|142.36
E48|  +0.43%  
B49|FLOKI
C49|https://coinranking.com/coin/fmHk13Rqw+floki-floki
D49|0.000270
E49|  +11.63%  
B50|VeChain
C50|https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet
D50|0.0360
E50|  +0.39%  
D51|39.24
E51|  +0.57%  
'@ -split "`r?`n"

foreach ($line in $updates) {
    if ([string]::IsNullOrEmpty($line)) { continue }
    $sep = $line.IndexOf('|')
    $ref = $line.Substring(0, $sep)
    $newValue = $line.Substring($sep + 1)

    $cell = $ws.Range($ref)

    # Plain decimal text such as "0.140", "602.11" or "0.0000251" would
    # otherwise be auto-converted by Excel into a Number (silently dropping
    # the significant trailing/leading zeros the source data relies on).
    # Force the cell to Text first so the literal digits survive, then put
    # the style back to Normal so we don't leave a stray text-format style
    # applied to the cell.
    $isNumericLooking = $newValue -match '^[+-]?\d+(\.\d+)?$'

    if ($isNumericLooking) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    } else {
        $cell.Value = $newValue
    }
}
